$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -20.92472190907698
$ws.Range("C2").Value = 2.483593965696561
$ws.Range("D2").Value = -20.92472190907698
$ws.Range("E2").Value = -20.92472190907698
$ws.Range("F2").Value = -20.92472190907698
$ws.Range("G2").Value = -20.92472190907698
$ws.Range("H2").Value = -20.92472190907698
$ws.Range("I2").Value = -20.92472190907698
$ws.Range("J2").Value = -20.92472190907698
$ws.Range("K2").Value = -20.92472190907698
$ws.Range("B3").Value = -20.92472190907698
$ws.Range("C3").Value = -20.92472190907698
$ws.Range("D3").Value = -20.92472190907698
$ws.Range("E3").Value = -20.92472190907698
$ws.Range("F3").Value = -20.92472190907698
$ws.Range("G3").Value = -20.92472190907698
$ws.Range("H3").Value = -20.92472190907698
$ws.Range("I3").Value = 2.375235221906761
$ws.Range("J3").Value = -20.92472190907698
$ws.Range("K3").Value = -20.92472190907698
$ws.Range("B4").Value = -20.92472190907698
$ws.Range("C4").Value = 2.084359882788016
$ws.Range("D4").Value = 2.890878038831076
$ws.Range("E4").Value = -20.92472190907698
$ws.Range("F4").Value = -20.92472190907698
$ws.Range("G4").Value = -20.92472190907698
$ws.Range("H4").Value = 1.84670166870533
$ws.Range("I4").Value = -20.92472190907698
$ws.Range("J4").Value = 2.379470215192558
$ws.Range("K4").Value = -20.92472190907698
$ws.Range("B5").Value = -20.92472190907698
$ws.Range("C5").Value = 1.012189513978905
$ws.Range("D5").Value = -20.92472190907698
$ws.Range("E5").Value = -20.92472190907698
$ws.Range("F5").Value = -20.92472190907698
$ws.Range("G5").Value = 2.112782985595762
$ws.Range("H5").Value = -20.92472190907698
$ws.Range("I5").Value = -20.92472190907698
$ws.Range("J5").Value = -20.92472190907698
$ws.Range("K5").Value = -20.92472190907698
$ws.Range("B6").Value = -20.92472190907698
$ws.Range("C6").Value = -20.92472190907698
$ws.Range("D6").Value = -20.92472190907698
$ws.Range("E6").Value = -20.92472190907698
$ws.Range("F6").Value = -20.92472190907698
$ws.Range("G6").Value = -20.92472190907698
$ws.Range("H6").Value = -20.92472190907698
$ws.Range("I6").Value = -20.92472190907698
$ws.Range("J6").Value = -20.92472190907698
$ws.Range("K6").Value = -20.92472190907698
$ws.Range("B7").Value = 2.990205230528271
$ws.Range("C7").Value = -20.92472190907698
$ws.Range("D7").Value = -20.92472190907698
$ws.Range("E7").Value = -20.92472190907698
$ws.Range("F7").Value = -20.92472190907698
$ws.Range("G7").Value = -20.92472190907698
$ws.Range("H7").Value = -20.92472190907698
$ws.Range("I7").Value = -20.92472190907698
$ws.Range("J7").Value = -20.92472190907698
$ws.Range("K7").Value = -20.92472190907698
$ws.Range("B8").Value = -20.92472190907698
$ws.Range("C8").Value = -20.92472190907698
$ws.Range("D8").Value = -20.92472190907698
$ws.Range("E8").Value = 2.901272636226032
$ws.Range("F8").Value = -20.92472190907698
$ws.Range("G8").Value = -20.92472190907698
$ws.Range("H8").Value = -20.92472190907698
$ws.Range("I8").Value = -20.92472190907698
$ws.Range("J8").Value = -20.92472190907698
$ws.Range("K8").Value = -20.92472190907698
$ws.Range("B9").Value = 3.591454514963458
$ws.Range("C9").Value = -20.92472190907698
$ws.Range("D9").Value = -20.92472190907698
$ws.Range("E9").Value = -20.92472190907698
$ws.Range("F9").Value = -20.92472190907698
$ws.Range("G9").Value = -20.92472190907698
$ws.Range("H9").Value = -20.92472190907698
$ws.Range("I9").Value = -20.92472190907698
$ws.Range("J9").Value = -20.92472190907698
$ws.Range("K9").Value = -20.92472190907698
$ws.Range("B10").Value = -20.92472190907698
$ws.Range("C10").Value = -20.92472190907698
$ws.Range("D10").Value = -20.92472190907698
$ws.Range("E10").Value = -20.92472190907698
$ws.Range("F10").Value = -20.92472190907698
$ws.Range("G10").Value = -20.92472190907698
$ws.Range("H10").Value = -20.92472190907698
$ws.Range("I10").Value = 1.568554287462386
$ws.Range("J10").Value = -20.92472190907698
$ws.Range("K10").Value = 2.23531152006105
$ws.Range("B11").Value = -20.92472190907698
$ws.Range("C11").Value = -20.92472190907698
$ws.Range("D11").Value = -20.92472190907698
$ws.Range("E11").Value = 1.980211986266016
$ws.Range("F11").Value = -20.92472190907698
$ws.Range("G11").Value = 2.606147908090585
$ws.Range("H11").Value = -20.92472190907698
$ws.Range("I11").Value = -20.92472190907698
$ws.Range("J11").Value = -20.92472190907698
$ws.Range("K11").Value = 1.35953768094531
$ws.Range("B12").Value = -20.92472190907698
$ws.Range("C12").Value = -20.92472190907698
$ws.Range("D12").Value = -20.92472190907698
$ws.Range("E12").Value = -20.92472190907698
$ws.Range("F12").Value = -20.92472190907698
$ws.Range("G12").Value = -20.92472190907698
$ws.Range("H12").Value = -20.92472190907698
$ws.Range("I12").Value = -20.92472190907698
$ws.Range("J12").Value = -20.92472190907698
$ws.Range("K12").Value = -20.92472190907698
$ws.Range("B13").Value = -20.92472190907698
$ws.Range("C13").Value = -20.92472190907698
$ws.Range("D13").Value = -20.92472190907698
$ws.Range("E13").Value = 1.65355794420917
$ws.Range("F13").Value = -20.92472190907698
$ws.Range("G13").Value = -20.92472190907698
$ws.Range("H13").Value = -20.92472190907698
$ws.Range("I13").Value = -20.92472190907698
$ws.Range("J13").Value = 2.266007795558398
$ws.Range("K13").Value = 1.663196808937122
$ws.Range("B14").Value = -20.92472190907698
$ws.Range("C14").Value = -20.92472190907698
$ws.Range("D14").Value = 1.603581719150747
$ws.Range("E14").Value = -20.92472190907698
$ws.Range("F14").Value = -20.92472190907698
$ws.Range("G14").Value = -20.92472190907698
$ws.Range("H14").Value = -20.92472190907698
$ws.Range("I14").Value = -20.92472190907698
$ws.Range("J14").Value = -20.92472190907698
$ws.Range("K14").Value = 2.112015313845144
$ws.Range("B15").Value = -20.92472190907698
$ws.Range("C15").Value = -20.92472190907698
$ws.Range("D15").Value = -0.009008076066901896
$ws.Range("E15").Value = -20.92472190907698
$ws.Range("F15").Value = -20.92472190907698
$ws.Range("G15").Value = -20.92472190907698
$ws.Range("H15").Value = -20.92472190907698
$ws.Range("I15").Value = -20.92472190907698
$ws.Range("J15").Value = -20.92472190907698
$ws.Range("K15").Value = -20.92472190907698
$ws.Range("B16").Value = -20.92472190907698
$ws.Range("C16").Value = -20.92472190907698
$ws.Range("D16").Value = -20.92472190907698
$ws.Range("E16").Value = -20.92472190907698
$ws.Range("F16").Value = -20.92472190907698
$ws.Range("G16").Value = -20.92472190907698
$ws.Range("H16").Value = -20.92472190907698
$ws.Range("I16").Value = -20.92472190907698
$ws.Range("J16").Value = 2.307034975097528
$ws.Range("K16").Value = -20.92472190907698
$ws.Range("B17").Value = -20.92472190907698
$ws.Range("C17").Value = 0.6843237819080921
$ws.Range("D17").Value = -0.1394761981728255
$ws.Range("E17").Value = -20.92472190907698
$ws.Range("F17").Value = -20.92472190907698
$ws.Range("G17").Value = -20.92472190907698
$ws.Range("H17").Value = 0.5661040011597075
$ws.Range("I17").Value = 0.9058215371558639
$ws.Range("J17").Value = 1.276209314140831
$ws.Range("K17").Value = -20.92472190907698
$ws.Range("B18").Value = -20.92472190907698
$ws.Range("C18").Value = -20.92472190907698
$ws.Range("D18").Value = -20.92472190907698
$ws.Range("E18").Value = -20.92472190907698
$ws.Range("F18").Value = -20.92472190907698
$ws.Range("G18").Value = -20.92472190907698
$ws.Range("H18").Value = 0.4715781744511908
$ws.Range("I18").Value = 0.9712826301348346
$ws.Range("J18").Value = 1.387343674294489
$ws.Range("K18").Value = -20.92472190907698
$ws.Range("B19").Value = -20.92472190907698
$ws.Range("C19").Value = -20.92472190907698
$ws.Range("D19").Value = 1.744710989833717
$ws.Range("E19").Value = -20.92472190907698
$ws.Range("F19").Value = -20.92472190907698
$ws.Range("G19").Value = -20.92472190907698
$ws.Range("H19").Value = 1.908839369168021
$ws.Range("I19").Value = 2.09291382538026
$ws.Range("J19").Value = -20.92472190907698
$ws.Range("K19").Value = -20.92472190907698
$ws.Range("B20").Value = -20.92472190907698
$ws.Range("C20").Value = 1.662359008377625
$ws.Range("D20").Value = 2.101282227170594
$ws.Range("E20").Value = -20.92472190907698
$ws.Range("F20").Value = 4.321927406347161
$ws.Range("G20").Value = -20.92472190907698
$ws.Range("H20").Value = 2.186908727523374
$ws.Range("I20").Value = 1.905132612373997
$ws.Range("J20").Value = -20.92472190907698
$ws.Range("K20").Value = 2.388243699350696
$ws.Range("B21").Value = -20.92472190907698
$ws.Range("C21").Value = 1.75573312505251
$ws.Range("D21").Value = -20.92472190907698
$ws.Range("E21").Value = 2.442955317195491
$ws.Range("F21").Value = -20.92472190907698
$ws.Range("G21").Value = 3.260926734770137
$ws.Range("H21").Value = 2.386182187737784
$ws.Range("I21").Value = -20.92472190907698
$ws.Range("J21").Value = -20.92472190907698
$ws.Range("K21").Value = -20.92472190907698

Write-Output "done"
